# Updates cryptos list prices / volume percentages (GitHub Actions style refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (preventing Excel's automatic
# numeric coercion of number-looking strings like "1.00" or "0.870"),
# then clear the temporary text-format style so the cell ends up with
# no style attribute, matching the original (unstyled) data cells.
function Set-TextValue($ws, $cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

$ws.Range("D2").Value = "54.437.24"
$ws.Range("E2").Value = "  -2.71%  "
$ws.Range("D3").Value = "2.286.13"
$ws.Range("E3").Value = "  -2.85%  "
Set-TextValue $ws "D4" "1.00"
$ws.Range("E4").Value = "  -0.06%  "
Set-TextValue $ws "D5" "493.74"
$ws.Range("E5").Value = "  -2.16%  "
Set-TextValue $ws "D6" "127.17"
$ws.Range("E6").Value = "  -2.46%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E8").Value = "  -1.77%  "
$ws.Range("D9").Value = "2.284.97"
$ws.Range("E9").Value = "  -3.39%  "
Set-TextValue $ws "D10" "0.0943"
$ws.Range("E10").Value = "  -3.02%  "
Set-TextValue $ws "D11" "0.150"
$ws.Range("E11").Value = "  +0.34%  "
Set-TextValue $ws "D12" "0.323"
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("E13").Value = "  -3.88%  "
$ws.Range("D14").Value = "2.687.08"
$ws.Range("E14").Value = "  -3.03%  "
Set-TextValue $ws "D15" "21.55"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").Value = "54.318.91"
$ws.Range("E16").Value = "  -2.84%  "
$ws.Range("E17").Value = "  -2.50%  "
$ws.Range("D18").Value = "2.289.16"
$ws.Range("E18").Value = "  -3.81%  "
Set-TextValue $ws "D19" "9.98"
$ws.Range("E19").Value = "  +0.69%  "
Set-TextValue $ws "D20" "4.06"
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws "D21" "6.50"
$ws.Range("E21").Value = "  +5.06%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws "D22" "303.53"
$ws.Range("E22").Value = "  -2.44%  "
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("E24").Value = "  -2.66%  "
$ws.Range("E25").Value = "  -2.88%  "
$ws.Range("E26").Value = "  +0.62%  "
Set-TextValue $ws "D27" "0.373"
$ws.Range("E27").Value = "  +0.82%  "
$ws.Range("D28").Value = "2.391.39"
$ws.Range("E28").Value = "  -3.00%  "
$ws.Range("E29").Value = "  +2.06%  "
Set-TextValue $ws "D30" "7.10"
$ws.Range("E30").Value = "  -0.51%  "
Set-TextValue $ws "D31" "168.99"
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("E32").Value = "  -2.56%  "
$ws.Range("D33").Value = "0.0₃0684"
$ws.Range("E33").Value = "  -3.12%  "
Set-TextValue $ws "D34" "5.89"
$ws.Range("E34").Value = "  +2.42%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("E37").Value = "  +0.99%  "
Set-TextValue $ws "D38" "17.59"
$ws.Range("E38").Value = "  -0.46%  "
Set-TextValue $ws "D39" "1.20"
$ws.Range("E39").Value = "  +1.92%  "
Set-TextValue $ws "D40" "0.870"
$ws.Range("E40").Value = "  +3.65%  "
Set-TextValue $ws "D41" "3.63"
$ws.Range("E41").Value = "  -0.59%  "
Set-TextValue $ws "D42" "35.56"
$ws.Range("E42").Value = "  -1.73%  "
$ws.Range("E43").Value = "  +0.89%  "
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("E45").Value = "  -0.15%  "
Set-TextValue $ws "D46" "128.48"
$ws.Range("E46").Value = "  +2.42%  "
Set-TextValue $ws "D47" "4.80"
$ws.Range("E47").Value = "  -1.14%  "
Set-TextValue $ws "D48" "0.0890"
$ws.Range("E48").Value = "  -0.65%  "
Set-TextValue $ws "D49" "0.543"
$ws.Range("E49").Value = "  -2.62%  "
Set-TextValue $ws "D50" "239.67"
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("E51").Value = "  -0.10%  "
